# Update the name in row 6 (Goli Kushwah -> Nitin Kushwah)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "Nitin Kushwah"

# Reflect the final active cell selection recorded in the sheet view (E12)
$ws.Range("E12").Select()
